$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing rows (33-35) that no longer exist in the refreshed model run
$ws.Rows("33:35").Delete()

# Update rows 2-32 with the refreshed feature-importance values (second model eval)
$data = @(
    ,@(2, "sstk_std", 0.04971417466192256, "Random Forest")
    ,@(3, "seq", 0.03857920496335165, "Random Forest")
    ,@(4, "teq", 0.03777455353823508, "Random Forest")
    ,@(5, "rest_sum_diff", 0.03096279781564133, "Random Forest")
    ,@(6, "xopr", 0.02940931022326146, "Random Forest")
    ,@(7, "ceq", 0.02590396875357646, "Random Forest")
    ,@(8, "caps", 0.02485518795836046, "Random Forest")
    ,@(9, "fopo_std", 0.02323817423979421, "Random Forest")
    ,@(10, "ceqt", 0.02185509492567017, "Random Forest")
    ,@(11, "pi_std", 0.02089165761873181, "Random Forest")
    ,@(12, "icapt", 0.02078584915138635, "Random Forest")
    ,@(13, "at", 0.02077248944233076, "Random Forest")
    ,@(14, "tstk", 0.01928175878266098, "Random Forest")
    ,@(15, "cogs", 0.01925306366437383, "Random Forest")
    ,@(16, "revt", 0.01893172555107008, "Random Forest")
    ,@(17, "sec_trt1m_std", 0.01886629464409779, "Random Forest")
    ,@(18, "rest_count", 0.01752467437758706, "Random Forest")
    ,@(19, "rest_a_count_of_diffs", 0.01700241946746915, "Random Forest")
    ,@(20, "lse", 0.01635374448785916, "Random Forest")
    ,@(21, "rest_count_of_diffs", 0.01623119717182976, "Random Forest")
    ,@(22, "cogs_std", 0.01579639654606069, "Random Forest")
    ,@(23, "st_per_growth", 0.01568318366536436, "Random Forest")
    ,@(24, "sstk", 0.01490504208014596, "Random Forest")
    ,@(25, "lct", 0.01431435953686472, "Random Forest")
    ,@(26, "gp", 0.01410872410330963, "Random Forest")
    ,@(27, "dilavx_std", 0.012630089981958, "Random Forest")
    ,@(28, "xopr_std", 0.01122832853933041, "Random Forest")
    ,@(29, "xsga", 0.01119838196101883, "Random Forest")
    ,@(30, "tstk_std", 0.01118596017342582, "Random Forest")
    ,@(31, "rect_std", 0.01060450737039614, "Random Forest")
    ,@(32, "rat_spcsrc", 0.001111, "Manual Addition")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
